# Update event stats on the "展览" (Exhibition) sheet and the
# "全部类型" (All types) sheet, which both list the same events.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

# Row 3: 南宁·2024良牙动漫秋季盛典（秋典）
$ws1.Range("F3").Value = 5605
$ws1.Range("G3").Value = 58

# Row 5: 广西·THO04-永夜廻想
$ws1.Range("F5").Value = 69

# Row 9: 南宁·万圣漫控嘉年华10
$ws1.Range("F9").Value = 531

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

# Row 3: 南宁·2024良牙动漫秋季盛典（秋典）
$ws4.Range("F3").Value = 5605
$ws4.Range("G3").Value = 58

# Row 6: 广西·THO04-永夜廻想
$ws4.Range("F6").Value = 69

# Row 11: 南宁·万圣漫控嘉年华10
$ws4.Range("F11").Value = 531
